# Generate Report for Handoff
# Updates the localization-status report with a new handoff id (GUID-like
# token) and refreshed timestamps for the latest handoff, replacing every
# occurrence of the old identifiers across all three worksheets.

$wb = $excel.ActiveWorkbook

$oldId = "124fb479-6146-435e-b80f-cc68695bdfd1"
$newId = "992a7d4b-62bc-4fe2-8532-89a06cdd18d3"

$oldZhHash = "df572b6ee3ef4879eddc31656c026b70666e15b0"
$newZhHash = "197741deac8efae800ae553eabccf86cb00e4387"

$oldDeHash = "df572b6ee3ef4879eddc31656c026b70666e15b0"
$newDeHash = "197741deac8efae800ae553eabccf86cb00e4387"

$oldMdName = "$oldId.md"
$newMdName = "$newId.md"

$oldZhXlf = "$oldId.$oldZhHash.zh-cn.xlf"
$newZhXlf = "$newId.$newZhHash.zh-cn.xlf"

$oldDeXlf = "$oldId.$oldDeHash.de-de.xlf"
$newDeXlf = "$newId.$newDeHash.de-de.xlf"

$oldOverviewDate = "2016-45-19 06:45:37"
$newOverviewDate = "2016-45-19 06:45:59"

$oldZhDate = "2016-03-19 06:45:34"
$newZhDate = "2016-03-19 06:45:56"

$oldDeDate = "2016-03-19 06:45:37"
$newDeDate = "2016-03-19 06:45:59"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newMdName
$wsOverview.Range("D2").Value = $newOverviewDate

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = $newZhDate
$wsZh.Hyperlinks.Item(1).TextToDisplay = $newMdName
$wsZh.Hyperlinks.Item(3).TextToDisplay = $newZhXlf

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = $newDeDate
$wsDe.Hyperlinks.Item(1).TextToDisplay = $newMdName
$wsDe.Hyperlinks.Item(3).TextToDisplay = $newDeXlf
